$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.409920036069533
    "C2" = 0.560620771757572
    "D2" = 0.534523033462853
    "E2" = 0.515083493490671
    "F2" = 0.749968017929011
    "G2" = 0.69441259681887
    "H2" = 0.684679095793903
    "I2" = 0.31598974774536
    "J2" = 0.312665570258618
    "K2" = 0.407749263806461
    "L2" = 0.555767853656013
    "M2" = 0.856324037927879
    "N2" = 0.12481239914949
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
